$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record ("Black Amber" ciruela) was added to the feed.
# In the source data it lands between the existing rows 150 and 151
# (original row order becomes 150, NEW, old-151, old-152, old-153, old-154),
# so insert a fresh row at 151 and push everything from the old row 151
# down by one.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record.
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 44610
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100103
$ws.Range("H151").Value = "Frutos de hueso (carozo)"
$ws.Range("I151").Value = 100103002
$ws.Range("J151").Value = "Ciruela"
$ws.Range("K151").Value = "Black Amber"
$ws.Range("L151").Value = "Segunda"
$ws.Range("M151").Value = 400
$ws.Range("N151").Value = 13000
$ws.Range("O151").Value = 13000
$ws.Range("P151").Value = 13000
$ws.Range("Q151").Value = "`$/caja 15 kilos granel"
$ws.Range("R151").Value = "Región de O'Higgins"
$ws.Range("S151").Value = 867
$ws.Range("T151").Value = 15
